$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 56: RMA/ship motor/encoder entry (1 hour, dated 3/8/2010)
# Copy formatting from the row above (A55) so the date cell keeps the
# same number format / style as the rest of column A, then set values.
$ws.Range("A55").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A56").Value = 40245
$ws.Range("B56").Value = 1
$ws.Range("C56").Value = "RMA/ship motor/encoder"

# Selection moves to the next empty row, same as after manual entry
$ws.Range("A57").Select()
